$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Value" column (B) for the "drill" and "kick" parameter rows:
# Row 13: drill / True  -> drill / drill
# Row 14: drill / False -> drill / None
# Row 15: kick  / True  -> kick  / kick
# Row 16: kick  / False -> kick  / None
$ws.Range("B13").Value = "drill"
$ws.Range("B14").Value = "None"
$ws.Range("B15").Value = "kick"
$ws.Range("B16").Value = "None"
